# Hoàn thiện Super Controller thiếu AddProject, EditProject chưa hoạt động,
# mới copy code sang PTLT
#
# - Rename "Sheet1" -> "All"
# - Column C changes from "labid" (numeric ids) to "subbid" (sub-project codes)
# - Add a new sheet "PTLT" (copy of the header row) after "All", and make it
#   the active sheet/tab

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "All"

# Column C: labid -> subbid
$ws1.Range("C1").Value = "subbid"
$ws1.Range("C2").Value = "PTLT1"
$ws1.Range("C3").Value = "PTCK1"
$ws1.Range("C4").Value = "PTTDH1"

$ws1.Range("A1:H1").Select() | Out-Null

# New sheet "PTLT" right after "All", seeded with the same header row
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "PTLT"

$ws1.Range("A1:H1").Copy() | Out-Null
$ws2.Range("A1").PasteSpecial() | Out-Null

$ws2.Range("D4").Select() | Out-Null
$ws2.Activate()
